$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-23 Wednesday" "2024-10-24 Thursday"

Replace-Text "83×65=5395" "74×67=4958"
Replace-Text "77×95=7315" "21×26=546"
Replace-Text "63×85=5355" "93×33=3069"
Replace-Text "86×76=6536" "66×88=5808"
Replace-Text "27×33=891" "18×17=306"

Replace-Text "92×19=1748" "80×99=7920"
Replace-Text "39×64=2496" "62×85=5270"
Replace-Text "62×52=3224" "85×81=6885"
Replace-Text "93×28=2604" "23×15=345"
Replace-Text "74×54=3996" "49×11=539"

Replace-Text "33×11=363" "96×79=7584"
Replace-Text "31×24=744" "52×85=4420"
Replace-Text "27×15=405" "95×67=6365"
Replace-Text "23×49=1127" "16×25=400"
Replace-Text "15×87=1305" "43×56=2408"

Replace-Text "35×73=2555" "88×31=2728"
Replace-Text "21×50=1050" "42×97=4074"
Replace-Text "65×64=4160" "94×49=4606"
Replace-Text "60×38=2280" "38×21=798"
Replace-Text "56×70=3920" "13×80=1040"

Replace-Text "66×39=2574" "80×51=4080"
Replace-Text "12×50=600" "40×16=640"
Replace-Text "67×99=6633" "45×28=1260"
Replace-Text "70×60=4200" "58×96=5568"
Replace-Text "90×42=3780" "27×44=1188"
